# Word COM-interop script implementing the commit:
#   "added new elements for mailing address"
#
# Changes:
#   1. Update the letter date from September 19, 2025 to September 21, 2025.
#   2. Split the recipient address paragraph "2916 Lamory Pl, Santa Clara CA 95051"
#      into two paragraphs: "2916 Lamory Pl" and "Santa Clara, CA 95051".
#   3. Remove the now-superfluous blank "No Spacing" paragraph that followed
#      "...Board of Directors" in the signature block.

$d = $word.ActiveDocument

# --- 1. Update the date -----------------------------------------------
# Assigning .Range.Text directly (rather than Find/Replace) keeps the
# xml:space="preserve" attribute that Word originally wrote on the run.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "September 19, 2025") {
        $p.Range.Text = "September 21, 2025"
        break
    }
}

# --- 2. Split the street address from the city/state/zip --------------
# Find the recipient address paragraph (the first occurrence, in the
# letterhead block) and split it into two paragraphs, preserving the
# paragraph/run formatting (Arial, 22pt) on the new line.
$found = $false
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "2916 Lamory Pl, Santa Clara CA 95051") {
        $p.Range.Text = "2916 Lamory Pl"
        $p.Range.InsertParagraphAfter()
        $newPara = $p.Next()
        $newPara.Range.Text = "Santa Clara, CA 95051"
        $found = $true
        break
    }
}

# --- 3. Remove the blank paragraph after "Board of Directors" ---------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Townhomes at Nuevo Homeowners Association Board of Directors") {
        $blank = $p.Next()
        $blank.Range.Delete()
        break
    }
}
